$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Price cells hold text (not numbers) in the source sheet, so numeric-looking
# price strings are written with a leading apostrophe to force them to stay text,
# matching Excel's normal quote-prefix behavior instead of auto-converting to a number.
$ws.Range("D2").Value = "89.191.58"
$ws.Range("E2").Value = "  -2.16%  "
$ws.Range("D3").Value = "3.096.22"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'212.85"
$ws.Range("E5").Value = "  -2.87%  "
$ws.Range("D6").Value = "'622.19"
$ws.Range("E6").Value = "  -1.98%  "
$ws.Range("E7").Value = "  -7.19%  "
$ws.Range("D8").Value = "'0.817"
$ws.Range("E8").Value = "  +14.86%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -2.93%  "
$ws.Range("D11").Value = "'0.622"
$ws.Range("E11").Value = "  +8.77%  "
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'0.0000242"
$ws.Range("E13").Value = "  -6.70%  "
$ws.Range("E14").Value = "  -1.65%  "
$ws.Range("D15").Value = "88.879.91"
$ws.Range("E15").Value = "  -2.05%  "
$ws.Range("D16").Value = "'32.30"
$ws.Range("E16").Value = "  -2.32%  "
$ws.Range("E17").Value = "  -3.02%  "
$ws.Range("D18").Value = "3.108.54"
$ws.Range("E18").Value = "  -3.00%  "
$ws.Range("D19").Value = "'3.40"
$ws.Range("E19").Value = "  +1.84%  "
$ws.Range("E20").Value = "  -6.46%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'423.91"
$ws.Range("E22").Value = "  -3.56%  "
$ws.Range("E23").Value = "  -3.24%  "
$ws.Range("D24").Value = "'4.95"
$ws.Range("E24").Value = "  -1.63%  "
$ws.Range("E25").Value = "  +5.42%  "
$ws.Range("D26").Value = "'11.91"
$ws.Range("E26").Value = "  +0.99%  "
$ws.Range("D27").Value = "'82.45"
$ws.Range("E27").Value = "  +2.19%  "
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E30").Value = "  +8.68%  "
$ws.Range("D31").Value = "'1.08"
$ws.Range("E31").Value = "  +7.75%  "
$ws.Range("E32").Value = "  -3.09%  "
$ws.Range("E34").Value = "  -10.88%  "
$ws.Range("E35").Value = "  -3.82%  "
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  +3.84%  "
$ws.Range("E40").Value = "  -0.46%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E43").Value = "  -2.28%  "
$ws.Range("E44").Value = "  -6.01%  "
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  +4.74%  "
$ws.Range("E47").Value = "  +13.39%  "
$ws.Range("E48").Value = "  -1.93%  "
$ws.Range("E49").Value = "  -6.56%  "
$ws.Range("E51").Value = "  -5.37%  "
